# Weekly update: a new price record (week of 2022-07-07, serial date 44749)
# is inserted as the newest entry right after the already up-to-date rows
# (rows 2-6), pushing all the older rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7:45 down to 8:46
$ws.Rows(7).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 44749
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108003
$ws.Range("J7").Value = "Maracuyá"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 34000
$ws.Range("O7").Value = 34000
$ws.Range("P7").Value = 34000
$ws.Range("Q7").Value = "$/caja 18 kilos"
$ws.Range("R7").Value = "Región de Arica y Parinacota"
$ws.Range("S7").Value = 1889
$ws.Range("T7").Value = 18
